# Swap the "Periodo Mora" / "Valor Mora" data between rows 16 and 17,
# and the "Periodo Mora" data between rows 18 and 19, on Hoja1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
# ($wb.ActiveSheet would also work here since the workbook has a single sheet)

# Row 16 / 17: swap Periodo Mora (E) text and Valor Mora (F) numeric value
$ws.Range("E16").Value = "1908"
$ws.Range("E17").Value = "1909"
$ws.Range("F16").Value = 9822
$ws.Range("F17").Value = 36834

# Row 18 / 19: swap Periodo Mora (E) text
$ws.Range("E18").Value = "1912"
$ws.Range("E19").Value = "2001"
